# KiCad FS error + updated BOM
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4: Schottky note wording update
$ws.Range("D4").Value = "Schottky; 1N4148 should work as well"

# P1-P4 notes: Arduino stackable header pack -> (long version)
$ws.Range("D6").Value = "Arduino stackable header pack (long version)"
$ws.Range("D7").Value = "Arduino stackable header pack (long version)"
$ws.Range("D8").Value = "Arduino stackable header pack (long version)"
$ws.Range("D9").Value = "Arduino stackable header pack (long version)"

# R1-R3: value gets a leading space, and a note is added
$ws.Range("B10").Value = " 5.1K"
$ws.Range("D10").Value = "any reasonably close value should work fine (i.e. 4.7K)"

# U2, U3: swap which logic family is the primary part vs. the alternative note
$ws.Range("B13").Value = "74HCT541"
$ws.Range("D13").Value = "74LS541 should work as well"

# U5, U6: remove the alternative-part note
$ws.Range("D15").Value = ""

# D16 stays blank but becomes an actual (empty) string cell
$ws.Range("D16").Value = ""

# Update the active cell selection saved with the sheet
$ws.Range("A2").Select()
